$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A23").Value = "second sublist"
$ws.Range("A24").Value = "second sublist"
$ws.Range("B24").Value = "node of second sublist"
$ws.Range("A25").Value = "third sublist"
$ws.Range("A26").Value = "third sublist"
$ws.Range("B26").Value = "special characters 1&2-%*_0 are embedded"

$ws.Range("A23:B26").Select() | Out-Null
